# Scheduled data refresh: update market/profit figures across the
# per-craft-job leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 450.75
$ws.Range("I4").Value = 450.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 450.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -336.75
$ws.Range("N4").Value = ""
$ws.Range("H42").Value = 615.63635
$ws.Range("I42").Value = 176.3
$ws.Range("K42").Value = 528.9000000000001
$ws.Range("M42").Value = -298.9000000000001
$ws.Range("H48").Value = 1416.25
$ws.Range("I48").Value = 443.33334
$ws.Range("K48").Value = 1330.00002
$ws.Range("M48").Value = -1038.00002
$ws.Range("H56").Value = 1416.25
$ws.Range("I56").Value = 443.33334
$ws.Range("K56").Value = 1330.00002
$ws.Range("M56").Value = -796.0000199999999
$ws.Range("H70").Value = 2648.2
$ws.Range("J70").Value = 2648.2
$ws.Range("L70").Value = 7944.599999999999
$ws.Range("N70").Value = -8484.599999999999
$ws.Range("H73").Value = 2648.2
$ws.Range("J73").Value = 2648.2
$ws.Range("L73").Value = 7944.599999999999
$ws.Range("N73").Value = -9816.599999999999
$ws.Range("H86").Value = 6652
$ws.Range("I86").Value = 6559
$ws.Range("J86").Value = 6810.1
$ws.Range("K86").Value = 6559
$ws.Range("L86").Value = 6810.1
$ws.Range("M86").Value = -5436
$ws.Range("N86").Value = -9056.1
$ws.Range("H89").Value = 6652
$ws.Range("I89").Value = 6559
$ws.Range("J89").Value = 6810.1
$ws.Range("K89").Value = 32795
$ws.Range("L89").Value = 34050.5
$ws.Range("M89").Value = -27179
$ws.Range("N89").Value = -45282.5
$ws.Range("H137").Value = 3672.3125
$ws.Range("I137").Value = 2363.2
$ws.Range("J137").Value = 5854.1665
$ws.Range("K137").Value = 7089.599999999999
$ws.Range("L137").Value = 17562.4995
$ws.Range("M137").Value = -4539.599999999999
$ws.Range("N137").Value = -22662.4995
$ws.Range("H138").Value = 1246.5
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 31317858
$ws.Range("I61").Value = 45457410
$ws.Range("J61").Value = 210849.2
$ws.Range("K61").Value = 45457410
$ws.Range("L61").Value = 210849.2
$ws.Range("M61").Value = -45457198
$ws.Range("N61").Value = -211273.2
$ws.Range("H74").Value = 8937813
$ws.Range("J74").Value = 21947
$ws.Range("L74").Value = 21947
$ws.Range("N74").Value = -23695
$ws.Range("H77").Value = 8937813
$ws.Range("J77").Value = 21947
$ws.Range("L77").Value = 109735
$ws.Range("N77").Value = -118471
$ws.Range("H132").Value = 4665.41
$ws.Range("I132").Value = 2480.8462
$ws.Range("K132").Value = 7442.5386
$ws.Range("M132").Value = -4912.5386
$ws.Range("H136").Value = 31317858
$ws.Range("I136").Value = 45457410
$ws.Range("J136").Value = 210849.2
$ws.Range("K136").Value = 136372230
$ws.Range("L136").Value = 632547.6000000001
$ws.Range("M136").Value = -136369680
$ws.Range("N136").Value = -637647.6000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1775.2
$ws.Range("I86").Value = 1753.8235
$ws.Range("K86").Value = 1753.8235
$ws.Range("M86").Value = -630.8235
$ws.Range("H89").Value = 1775.2
$ws.Range("I89").Value = 1753.8235
$ws.Range("K89").Value = 8769.1175
$ws.Range("M89").Value = -3153.1175
$ws.Range("H105").Value = 2381.1428
$ws.Range("I105").Value = 2228
$ws.Range("J105").Value = 2496
$ws.Range("K105").Value = 2228
$ws.Range("L105").Value = 2496
$ws.Range("M105").Value = -481
$ws.Range("N105").Value = -5990

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1543.8846
$ws.Range("I7").Value = 129.78572
$ws.Range("K7").Value = 129.78572
$ws.Range("M7").Value = -16.78572
$ws.Range("H105").Value = 2500.6
$ws.Range("I105").Value = 2145.1428
$ws.Range("J105").Value = 3330
$ws.Range("K105").Value = 2145.1428
$ws.Range("L105").Value = 3330
$ws.Range("M105").Value = -398.1428000000001
$ws.Range("N105").Value = -6824
$ws.Range("H107").Value = 1269.8125
$ws.Range("I107").Value = 1128.1111
$ws.Range("J107").Value = 1452
$ws.Range("K107").Value = 1128.1111
$ws.Range("L107").Value = 1452
$ws.Range("M107").Value = 791.8888999999999
$ws.Range("N107").Value = -5292
$ws.Range("H114").Value = 76881.25
$ws.Range("J114").Value = 76881.25
$ws.Range("L114").Value = 76881.25
$ws.Range("N114").Value = -85559.25
$ws.Range("H132").Value = 4627.857
$ws.Range("I132").Value = 4580.1665
$ws.Range("K132").Value = 13740.4995
$ws.Range("M132").Value = -11210.4995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 229.14285
$ws.Range("I14").Value = 229.14285
$ws.Range("K14").Value = 687.4285500000001
$ws.Range("M14").Value = -514.4285500000001
$ws.Range("H44").Value = 56946
$ws.Range("I44").Value = 113695.664
$ws.Range("J44").Value = 196.33333
$ws.Range("K44").Value = 341086.992
$ws.Range("L44").Value = 588.99999
$ws.Range("M44").Value = -340688.992
$ws.Range("N44").Value = -1384.99999
$ws.Range("H132").Value = 2489.4443
$ws.Range("J132").Value = 2797
$ws.Range("L132").Value = 25173
$ws.Range("N132").Value = -30233

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 74352
$ws.Range("I44").Value = 74514
$ws.Range("K44").Value = 74514
$ws.Range("M44").Value = -73918
$ws.Range("H70").Value = 7308.8335
$ws.Range("I70").Value = 6049.375
$ws.Range("K70").Value = 6049.375
$ws.Range("M70").Value = -5779.375
$ws.Range("H73").Value = 7308.8335
$ws.Range("I73").Value = 6049.375
$ws.Range("K73").Value = 6049.375
$ws.Range("M73").Value = -5113.375
$ws.Range("H97").Value = 858.1818
$ws.Range("J97").Value = 1211.2858
$ws.Range("L97").Value = 1211.2858
$ws.Range("N97").Value = -2203.2858
$ws.Range("H113").Value = 4103.7
$ws.Range("I113").Value = 3012
$ws.Range("J113").Value = 4225
$ws.Range("K113").Value = 3012
$ws.Range("L113").Value = 4225
$ws.Range("M113").Value = -842
$ws.Range("N113").Value = -8565
$ws.Range("H132").Value = 33336178
$ws.Range("I132").Value = 38464172
$ws.Range("K132").Value = 115392516
$ws.Range("M132").Value = -115389986

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1842.5714
$ws.Range("I82").Value = 1700
$ws.Range("K82").Value = 1700
$ws.Range("M82").Value = -1339
$ws.Range("H85").Value = 1842.5714
$ws.Range("I85").Value = 1700
$ws.Range("K85").Value = 1700
$ws.Range("M85").Value = -452
$ws.Range("H93").Value = 31250974
$ws.Range("J93").Value = 998.5714
$ws.Range("L93").Value = 998.5714
$ws.Range("N93").Value = -3494.5714

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3325.1
$ws.Range("I132").Value = 3306.7778
$ws.Range("J132").Value = 3490
$ws.Range("K132").Value = 9920.3334
$ws.Range("L132").Value = 10470
$ws.Range("M132").Value = -7390.3334
$ws.Range("N132").Value = -15530
